$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain plain text,
# matching the source workbook where these are stored as inline strings
# (not numbers/percentages), then write the refreshed values scraped by
# the GitHub Actions job.
$textCells = @(
    "D2"
    "E2"
    "D3"
    "E3"
    "D4"
    "E4"
    "D5"
    "E5"
    "D6"
    "E6"
    "D7"
    "E7"
    "D8"
    "E8"
    "D9"
    "E9"
    "D10"
    "E10"
    "D11"
    "E11"
    "D12"
    "E12"
    "D13"
    "E13"
    "D14"
    "E14"
    "D15"
    "E15"
    "D16"
    "E16"
    "D17"
    "E17"
    "E18"
    "E19"
    "D20"
    "E20"
    "D21"
    "E21"
    "D22"
    "E22"
    "D23"
    "E23"
    "E24"
    "D25"
    "E25"
    "E26"
    "D27"
    "E27"
    "D39"
    "E39"
    "D40"
    "E40"
    "D41"
    "E41"
    "E42"
    "E43"
    "D44"
    "E44"
    "D45"
    "E45"
    "D46"
    "E46"
    "D47"
    "E47"
    "D48"
    "E48"
    "D49"
    "E49"
    "D50"
    "E50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "298.26"
$ws.Range("E2").Value = "0.08%"
$ws.Range("D3").Value = "31.34"
$ws.Range("E3").Value = "-0.03%"
$ws.Range("D4").Value = "5.104"
$ws.Range("E4").Value = "-0.26%"
$ws.Range("D5").Value = "0.08025"
$ws.Range("E5").Value = "9.56%"
$ws.Range("D6").Value = "2.444"
$ws.Range("E6").Value = "28.03%"
$ws.Range("D7").Value = "7.821"
$ws.Range("E7").Value = "0.80%"
$ws.Range("D8").Value = "3.804"
$ws.Range("E8").Value = "1.87%"
$ws.Range("D9").Value = "0.9201"
$ws.Range("E9").Value = "-0.64%"
$ws.Range("D10").Value = "0.1730"
$ws.Range("E10").Value = "3.31%"
$ws.Range("D11").Value = "0.07299"
$ws.Range("E11").Value = "3.68%"
$ws.Range("D12").Value = "0.08514"
$ws.Range("E12").Value = "7.15%"
$ws.Range("D13").Value = "0.03037"
$ws.Range("E13").Value = "1.11%"
$ws.Range("D14").Value = "0.09970"
$ws.Range("E14").Value = "0.42%"
$ws.Range("D15").Value = "0.001506"
$ws.Range("E15").Value = "0.95%"
$ws.Range("D16").Value = "0.006016"
$ws.Range("E16").Value = "-1.53%"
$ws.Range("D17").Value = "3.511"
$ws.Range("E17").Value = "1.64%"
$ws.Range("E18").Value = "0.92%"
$ws.Range("E19").Value = "1.79%"
$ws.Range("D20").Value = "0.1338"
$ws.Range("E20").Value = "1.81%"
$ws.Range("D21").Value = "4.624"
$ws.Range("E21").Value = "1.57%"
$ws.Range("D22").Value = "0.1618"
$ws.Range("E22").Value = "2.32%"
$ws.Range("D23").Value = "0.04632"
$ws.Range("E23").Value = "-0.25%"
$ws.Range("E24").Value = "2.69%"
$ws.Range("D25").Value = "0.004432"
$ws.Range("E25").Value = "-6.35%"
$ws.Range("E26").Value = "-7.46%"
$ws.Range("D27").Value = "0.0003430"
$ws.Range("E27").Value = "83.12%"
$ws.Range("D39").Value = "0.01795"
$ws.Range("E39").Value = "4.29%"
$ws.Range("D40").Value = "0.04486"
$ws.Range("E40").Value = "0.27%"
$ws.Range("D41").Value = "0.007021"
$ws.Range("E41").Value = "-1.62%"
$ws.Range("E42").Value = "0.69%"
$ws.Range("E43").Value = "1.61%"
$ws.Range("D44").Value = "0.009845"
$ws.Range("E44").Value = "-7.50%"
$ws.Range("D45").Value = "0.00006597"
$ws.Range("E45").Value = "5.88%"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.09%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "0.005238"
$ws.Range("E47").Value = "-48.72%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.8206"
$ws.Range("E48").Value = "11.05%"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").Value = "0.09%"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").Value = "0.16%"
